$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    if ($text -like "Modelo FlexBox*") {
        $p.Range.Font.HighlightColorIndex = 7
    } elseif ($text -like "Modelo Grid*") {
        $p.Range.Font.HighlightColorIndex = 3
    } elseif ($text -like "Media Queries. Responsive Design*") {
        $p.Range.Font.HighlightColorIndex = 7
    }
}
